$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 42 / 43: coin rows swapped (RenderToken <-> TheSandbox) ---
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.397"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.01%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5358"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.73%  "

# --- Price (D) / Volume(1h) (E) updates for all other rows ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.380.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.76%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.839.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.28%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.014"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.29%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.00%  "

$ws.Range("E6").Value = "  +1.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4742"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.56%  "

$ws.Range("E8").Value = "  +0.34%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07466"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.47%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8858"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.68%  "

$ws.Range("E11").Value = "  +0.51%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.836.96"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07373"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.481"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.93%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.79%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.583"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.03%  "

$ws.Range("E17").Value = "  +1.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008845"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.67%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.013"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.403.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.351"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.56%  "

$ws.Range("E23").Value = "  +1.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.071.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.914"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.96%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.26%  "

$ws.Range("E27").Value = "  +1.62%  "

$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.262"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.83%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08969"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7602"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.180"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.558"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.20%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.938"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.013"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.25%  "

$ws.Range("E37").Value = "  +1.78%  "

$ws.Range("E38").Value = "  +1.69%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01964"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.32%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.288"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.39%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1664"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.37%  "

$ws.Range("E45").Value = "  +1.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4980"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.77%  "

$ws.Range("E48").Value = "  +1.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.35%  "

$ws.Range("E50").Value = "  +0.55%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06320"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.32%  "
